$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
Write-Output ($tcs | Get-Member | Out-String)
